# Fruta / hortaliza, semanal
# Re-shuffle the weekly rows (2-8, 10-13; row 9 unchanged) across columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg).
#
# Mapping is: new row <- data currently sitting in old row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 11
    6  = 2
    7  = 12
    8  = 3
    10 = 13
    11 = 4
    12 = 10
    13 = 8
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the affected columns/rows before
# writing anything, since several destinations read from rows that will
# also be overwritten.
$snapshot = @{}
foreach ($row in $mapping.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$row").Value2
        }
        $snapshot[$row] = $rowData
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $src[$col]
    }
}
